$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 5 entirely (worst-fit-algorithm / seed 7914 run no longer present)
$ws.Rows.Item(5).Delete()

# Row 2 updates
$ws.Range("B2").Value = 81
$ws.Range("D2").Value = 67
$ws.Range("E2").Value = 87
$ws.Range("F2").Value = 0.7701149425287356
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1143
$ws.Range("K2").Value = 1056
$ws.Range("L2").Value = 87
$ws.Range("M2").Value = 6
$ws.Range("N2").Value = 4
$ws.Range("Q2").Value = 49.34388446807861
$ws.Range("R2").Value = 0.3690639998689534
$ws.Range("S2").Value = 0.2017369417862839
$ws.Range("T2").Value = 0.3
$ws.Range("U2").Value = 0.3333333333333333

# Row 3 updates
$ws.Range("B3").Value = 81
$ws.Range("D3").Value = 67
$ws.Range("E3").Value = 103
$ws.Range("F3").Value = 0.6504854368932039
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("J3").Value = 1143
$ws.Range("K3").Value = 1040
$ws.Range("L3").Value = 103
$ws.Range("M3").Value = 9
$ws.Range("N3").Value = 4
$ws.Range("Q3").Value = 53.64048075675964
$ws.Range("R3").Value = 0.357070501514946
$ws.Range("S3").Value = 0.1867144540455617
$ws.Range("T3").Value = 0.45
$ws.Range("U3").Value = 0.3333333333333333

# Row 4 updates
$ws.Range("B4").Value = 81
$ws.Range("C4").Value = "nord-algorithm"
$ws.Range("D4").Value = 67
$ws.Range("E4").Value = 103
$ws.Range("F4").Value = 0.6504854368932039
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("J4").Value = 1143
$ws.Range("K4").Value = 1040
$ws.Range("L4").Value = 103
$ws.Range("M4").Value = 7
$ws.Range("N4").Value = 4
$ws.Range("Q4").Value = 52.26946973800659
$ws.Range("R4").Value = 0.4158672015814873
$ws.Range("S4").Value = 0.1872084415976796
$ws.Range("T4").Value = 0.35
$ws.Range("U4").Value = 0.3333333333333333
